# Refactor synthetic array /3: re-theme the "statut" color swatches from
# black/orange/green/red squares to blue/orange/green/red book emoji, and
# rename the "noir" label to "bleu" to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.UsedRange

foreach ($cell in $cells) {
    $v = $cell.Value2
    if ($null -ne $v -and $v.GetType().Name -eq "String") {
        if ($v -eq "⬛") {
            $cell.Value2 = "📘"
        } elseif ($v -eq "🟧") {
            $cell.Value2 = "📙"
        } elseif ($v -eq "🟩") {
            $cell.Value2 = "📗"
        } elseif ($v -eq "🟥") {
            $cell.Value2 = "📕"
        } elseif ($v -eq "noir") {
            $cell.Value2 = "bleu"
        }
    }
}
